$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 updates
$ws.Range("G4").Value = 2.1
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("AH4").Value = 17
$ws.Range("AP4").Value = 29
$ws.Range("AT4").Value = 2.25

# Row 11 updates
$ws.Range("Q11").Value = 2.15
$ws.Range("R11").Value = 1.67

# Row 12 updates
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("Q12").Value = 2.08
$ws.Range("R12").Value = 1.73
